$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 43908.6
$ws.Range("J3").Value = 43908.6
$ws.Range("L3").Value = 43908.6
$ws.Range("N3").Value = -44136.6
$ws.Range("H17").Value = 1205.4897
$ws.Range("I17").Value = 960
$ws.Range("J17").Value = 1210.6041
$ws.Range("K17").Value = 2880
$ws.Range("L17").Value = 3631.8123
$ws.Range("N17").Value = -3967.8123
$ws.Range("M17").Value = -2712
$ws.Range("H80").Value = 25641984
$ws.Range("I80").Value = 37037900
$ws.Range("J80").Value = 1175
$ws.Range("K80").Value = 111113700
$ws.Range("L80").Value = 3525
$ws.Range("M80").Value = -111112702
$ws.Range("N80").Value = -5521
$ws.Range("H83").Value = 25641984
$ws.Range("I83").Value = 37037900
$ws.Range("J83").Value = 1175
$ws.Range("K83").Value = 333341100
$ws.Range("L83").Value = 10575
$ws.Range("M83").Value = -333336108
$ws.Range("N83").Value = -20559
$ws.Range("H86").Value = 7694534
$ws.Range("I86").Value = 7694534
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 7694534
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = -7693411
$ws.Range("H89").Value = 7694534
$ws.Range("I89").Value = 7694534
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 38472670
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = -38467054
$ws.Range("H102").Value = 43908.6
$ws.Range("J102").Value = 43908.6
$ws.Range("L102").Value = 43908.6
$ws.Range("N102").Value = -50398.6
$ws.Range("H108").Value = 31260
$ws.Range("J108").Value = 31260
$ws.Range("L108").Value = 31260
$ws.Range("N108").Value = -38940
$ws.Range("H109").Value = 28276.666
$ws.Range("J109").Value = 28276.666
$ws.Range("L109").Value = 28276.666
$ws.Range("N109").Value = -31050.666
$ws.Range("H120").Value = 49714
$ws.Range("J120").Value = 49714
$ws.Range("L120").Value = 49714
$ws.Range("N120").Value = -59390
$ws.Range("H128").Value = 44170.5
$ws.Range("J128").Value = 44170.5
$ws.Range("L128").Value = 44170.5
$ws.Range("N128").Value = -54130.5
$ws.Range("H130").Value = 45601.6
$ws.Range("J130").Value = 45601.6
$ws.Range("L130").Value = 45601.6
$ws.Range("N130").Value = -55641.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29571.34
$ws.Range("I32").Value = 29407.77
$ws.Range("J32").Value = 31697.75
$ws.Range("K32").Value = 29407.77
$ws.Range("L32").Value = 31697.75
$ws.Range("M32").Value = -29120.77
$ws.Range("N32").Value = -32271.75
$ws.Range("H45").Value = 1755.409
$ws.Range("I45").Value = 1478.1111
$ws.Range("J45").Value = 3003.25
$ws.Range("K45").Value = 1478.1111
$ws.Range("L45").Value = 3003.25
$ws.Range("M45").Value = -1101.1111
$ws.Range("N45").Value = -3757.25
$ws.Range("H61").Value = 3141.6775
$ws.Range("I61").Value = 1742.5714
$ws.Range("K61").Value = 1742.5714
$ws.Range("M61").Value = -1530.5714
$ws.Range("H98").Value = 43351
$ws.Range("J98").Value = 43351
$ws.Range("L98").Value = 43351
$ws.Range("N98").Value = -49341
$ws.Range("H101").Value = 44096
$ws.Range("J101").Value = 44096
$ws.Range("L101").Value = 44096
$ws.Range("N101").Value = -50586
$ws.Range("H103").Value = 37332.57
$ws.Range("J103").Value = 37332.57
$ws.Range("L103").Value = 37332.57
$ws.Range("N103").Value = -39676.57
$ws.Range("H117").Value = 46523.715
$ws.Range("J117").Value = 46523.715
$ws.Range("L117").Value = 46523.715
$ws.Range("N117").Value = -55701.715
$ws.Range("H121").Value = 41061.75
$ws.Range("J121").Value = 41061.75
$ws.Range("L121").Value = 41061.75
$ws.Range("N121").Value = -44555.75
$ws.Range("H131").Value = 41006
$ws.Range("J131").Value = 41006
$ws.Range("L131").Value = 41006
$ws.Range("N131").Value = -51086
$ws.Range("H132").Value = 21740752
$ws.Range("I132").Value = 27778884
$ws.Range("K132").Value = 83336652
$ws.Range("M132").Value = -83334122
$ws.Range("H136").Value = 3141.6775
$ws.Range("I136").Value = 1742.5714
$ws.Range("K136").Value = 5227.7142
$ws.Range("M136").Value = -2677.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 39840.4
$ws.Range("J106").Value = 39840.4
$ws.Range("L106").Value = 39840.4
$ws.Range("N106").Value = -42364.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 964
$ws.Range("J16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("N16").Value = -1574
$ws.Range("H96").Value = 39104.46
$ws.Range("J96").Value = 39104.46
$ws.Range("L96").Value = 39104.46
$ws.Range("N96").Value = -44596.46
$ws.Range("H106").Value = 39313.6
$ws.Range("J106").Value = 39313.6
$ws.Range("L106").Value = 39313.6
$ws.Range("N106").Value = -41837.6
$ws.Range("H110").Value = 42195.5
$ws.Range("J110").Value = 42195.5
$ws.Range("L110").Value = 42195.5
$ws.Range("N110").Value = -50375.5
$ws.Range("H113").Value = 964
$ws.Range("J113").Value = 1000
$ws.Range("L113").Value = 1000
$ws.Range("N113").Value = -5340
$ws.Range("H132").Value = 51399.965
$ws.Range("I132").Value = 2159.5293
$ws.Range("J132").Value = 121157.25
$ws.Range("K132").Value = 6478.5879
$ws.Range("L132").Value = 363471.75
$ws.Range("M132").Value = -3948.5879
$ws.Range("N132").Value = -368531.75
$ws.Range("H134").Value = 90410
$ws.Range("I134").Value = 1603.8889
$ws.Range("J134").Value = 204589.28
$ws.Range("K134").Value = 4811.6667
$ws.Range("L134").Value = 613767.84
$ws.Range("M134").Value = -2276.6667
$ws.Range("N134").Value = -618837.84

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 642.8570999999999
$ws.Range("J4").Value = 1800
$ws.Range("L4").Value = 5400
$ws.Range("N4").Value = -5624
$ws.Range("H113").Value = 3407.75
$ws.Range("I113").Value = 8307.846
$ws.Range("J113").Value = 638.13043
$ws.Range("K113").Value = 24923.538
$ws.Range("L113").Value = 1914.39129
$ws.Range("M113").Value = -22753.538
$ws.Range("N113").Value = -6254.39129

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 500
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H130").Value = 45905.375
$ws.Range("J130").Value = 45905.375
$ws.Range("L130").Value = 45905.375
$ws.Range("N130").Value = -55945.375
$ws.Range("H132").Value = 4355.8965
$ws.Range("I132").Value = 1623
$ws.Range("J132").Value = 6576.375
$ws.Range("K132").Value = 4869
$ws.Range("L132").Value = 19729.125
$ws.Range("M132").Value = -2339
$ws.Range("N132").Value = -24789.125
$ws.Range("H138").Value = 49209.09
$ws.Range("J138").Value = 49209.09
$ws.Range("L138").Value = 49209.09
$ws.Range("N138").Value = -59489.09

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 40748
$ws.Range("J108").Value = 40748
$ws.Range("L108").Value = 40748
$ws.Range("N108").Value = -48428
$ws.Range("H109").Value = 35277
$ws.Range("J109").Value = 35277
$ws.Range("L109").Value = 35277
$ws.Range("N109").Value = -38051
$ws.Range("H110").Value = 32357.334
$ws.Range("J110").Value = 32357.334
$ws.Range("L110").Value = 32357.334
$ws.Range("N110").Value = -40537.334
$ws.Range("H118").Value = 38865.832
$ws.Range("J118").Value = 38865.832
$ws.Range("L118").Value = 38865.832
$ws.Range("N118").Value = -42179.832
$ws.Range("H132").Value = 3241.0881
$ws.Range("I132").Value = 1186.2858
$ws.Range("J132").Value = 4679.45
$ws.Range("K132").Value = 3558.8574
$ws.Range("L132").Value = 14038.35
$ws.Range("M132").Value = -1028.8574
$ws.Range("N132").Value = -19098.35
$ws.Range("H134").Value = 45017.54
$ws.Range("J134").Value = 45017.54
$ws.Range("L134").Value = 45017.54
$ws.Range("N134").Value = -55157.54

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 1308607.8
$ws.Range("I5").Value = 1334322.2
$ws.Range("J5").Value = 1250750
$ws.Range("K5").Value = 1334322.2
$ws.Range("L5").Value = 1250750
$ws.Range("M5").Value = -1334210.2
$ws.Range("N5").Value = -1250974
$ws.Range("H27").Value = 36980
$ws.Range("J27").Value = 36980
$ws.Range("L27").Value = 36980
$ws.Range("N27").Value = -37118
$ws.Range("H105").Value = 49607
$ws.Range("J105").Value = 49607
$ws.Range("L105").Value = 49607
$ws.Range("N105").Value = -56595
$ws.Range("H108").Value = 29695.2
$ws.Range("J108").Value = 29695.2
$ws.Range("L108").Value = 29695.2
$ws.Range("N108").Value = -37375.2
$ws.Range("H112").Value = 36674.2
$ws.Range("J112").Value = 36674.2
$ws.Range("L112").Value = 36674.2
$ws.Range("N112").Value = -39628.2
$ws.Range("H115").Value = 30451.572
$ws.Range("J115").Value = 30451.572
$ws.Range("L115").Value = 30451.572
$ws.Range("N115").Value = -33585.572
$ws.Range("H120").Value = 35210
$ws.Range("J120").Value = 35210
$ws.Range("L120").Value = 35210
$ws.Range("N120").Value = -44886
$ws.Range("H133").Value = 102728.5
$ws.Range("J133").Value = 102728.5
$ws.Range("L133").Value = 102728.5
$ws.Range("N133").Value = -112848.5
$ws.Range("H137").Value = 48999.668
$ws.Range("J137").Value = 48999.668
$ws.Range("L137").Value = 48999.668
$ws.Range("N137").Value = -59199.668
